$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K column (col G) values after recalculating std/mean and s_vals
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
